# Applies the crypto price-table refresh described by the commit:
# "Updated symbol list on Wed Dec 28 06:57:09 UTC 2022 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = '''243.05'
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.Value = '''23.52'
$cell.Style = "Normal"
$cell = $ws.Range("D4")
$cell.Value = '''5.299'
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.Value = '''6.475'
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.Value = '''3.336'
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.Value = '''0.8100'
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.Value = '''0.8832'
$cell.Style = "Normal"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$cell = $ws.Range("D10")
$cell.Value = '''0.1378'
$cell.Style = "Normal"
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$cell = $ws.Range("D11")
$cell.Value = '''0.07289'
$cell.Style = "Normal"
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$cell = $ws.Range("D12")
$cell.Value = '''0.03085'
$cell.Style = "Normal"
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$cell = $ws.Range("D13")
$cell.Value = '''0.03060'
$cell.Style = "Normal"
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$cell = $ws.Range("D14")
$cell.Value = '''0.09320'
$cell.Style = "Normal"
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$cell = $ws.Range("D15")
$cell.Value = '''3.866'
$cell.Style = "Normal"
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$cell = $ws.Range("D16")
$cell.Value = '''0.001543'
$cell.Style = "Normal"
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$cell = $ws.Range("D17")
$cell.Value = '''0.04699'
$cell.Style = "Normal"
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$cell = $ws.Range("D18")
$cell.Value = '''0.0006046'
$cell.Style = "Normal"
$ws.Range("E18").Value = '17OneONE'
$cell = $ws.Range("D19")
$cell.Value = '''0.006014'
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.Value = '''0.001304'
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.Value = '''0.004599'
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.Value = '''0.00008808'
$cell.Style = "Normal"
$ws.Range("E22").Value = '21NitroExNTXBestin24h'
$cell = $ws.Range("D23")
$cell.Value = '''3.577'
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.Value = '''2.143'
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.Value = '''0.3181'
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.Value = '''0.1320'
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.Value = '''0.03766'
$cell.Style = "Normal"
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$cell = $ws.Range("D41")
$cell.Value = '''0.1051'
$cell.Style = "Normal"
$ws.Range("E41").Value = '40BKEXTokenBKK'
$cell = $ws.Range("D42")
$cell.Value = '''0.002562'
$cell.Style = "Normal"
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$cell = $ws.Range("D43")
$cell.Value = '''0.003180'
$cell.Style = "Normal"
$ws.Range("E43").Value = '42KickTokenKICKWorstin24h'
$cell = $ws.Range("D44")
$cell.Value = '''0.007132'
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.Value = '''0.00005489'
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.Value = '''0.6006'
$cell.Style = "Normal"
$ws.Range("E48").Value = '47BOLOBOLO'
